# Leave Card update - 4/18/2023 4:57 PM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Clear the stale BALANCE formula results for the 2018-2023 entries
# (row 10 is the first "2018" section row, rows 11-134 are the data rows)
$ws.Range("E10:E134").ClearContents()

# Update the printed page scale
$ws.PageSetup.Zoom = 88

# Update certifying officer in the footer
$ws.PageSetup.CenterFooter = "`nCERTIFIED CORRECT BY: &UNANETTE B. SUSA&U`n                                           OIC-HRMO"

# Restore the last-used selection/view position
$ws.Application.Goto($ws.Range("D15"))
